$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the article/status table (header + 6 data rows, columns A:D).
$data = @(
    @("Artikelnummer", "Menge", "Empfänger", "Status"),
    @("KL5011", 2, "PL", 0),
    @("KL5012", 3, "HUN", 0),
    @("KL5013", 4, "HUN", 0),
    @("KL5014", 5, "HUN", 1),
    @("KL5015", 6, "3 HUN; 3 PL", 0),
    @("KL5016", 7, "PL", 0)
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $rowValues[$c]
    }
}

# Widen the "Empfänger" column (C) so the longest entry fits.
$ws.Columns.Item(3).ColumnWidth = 15

# Match the selection/active range saved with the populated sheet.
$ws.Range("A1:D7").Select() | Out-Null
